$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3430
$ws.Range("I62").Value = 2600
$ws.Range("J62").Value = 3637.5
$ws.Range("K62").Value = 2600
$ws.Range("L62").Value = 3637.5
$ws.Range("M62").Value = -1976
$ws.Range("N62").Value = -4885.5
$ws.Range("H65").Value = 3430
$ws.Range("I65").Value = 2600
$ws.Range("J65").Value = 3637.5
$ws.Range("K65").Value = 13000
$ws.Range("L65").Value = 18187.5
$ws.Range("M65").Value = -9880
$ws.Range("N65").Value = -24427.5
$ws.Range("H129").Value = 718.8
$ws.Range("I129").Value = 340.6
$ws.Range("K129").Value = 1021.8
$ws.Range("M129").Value = 3978.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3136.945
$ws.Range("I32").Value = 2536.2278
$ws.Range("J32").Value = 7091.6665
$ws.Range("K32").Value = 2536.2278
$ws.Range("L32").Value = 7091.6665
$ws.Range("M32").Value = -2249.2278
$ws.Range("N32").Value = -7665.6665
$ws.Range("H45").Value = 1399
$ws.Range("I45").Value = 1252.0714
$ws.Range("J45").Value = 1692.8572
$ws.Range("K45").Value = 1252.0714
$ws.Range("L45").Value = 1692.8572
$ws.Range("M45").Value = -875.0714
$ws.Range("N45").Value = -2446.8572
$ws.Range("H74").Value = 1157.4324
$ws.Range("I74").Value = 782.4828
$ws.Range("K74").Value = 782.4828
$ws.Range("M74").Value = 91.5172
$ws.Range("H77").Value = 1157.4324
$ws.Range("I77").Value = 782.4828
$ws.Range("K77").Value = 3912.414
$ws.Range("M77").Value = 455.586
$ws.Range("H110").Value = 1212.6666
$ws.Range("I110").Value = 819
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 819
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 1226
$ws.Range("N110").Value = -6090
$ws.Range("H122").Value = 1645.3055
$ws.Range("I122").Value = 1112.5769
$ws.Range("J122").Value = 3030.4
$ws.Range("K122").Value = 3337.7307
$ws.Range("L122").Value = 9091.2
$ws.Range("M122").Value = -887.7307000000001
$ws.Range("N122").Value = -13991.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1366.3334
$ws.Range("I20").Value = 1010.5833
$ws.Range("J20").Value = 2315
$ws.Range("K20").Value = 1010.5833
$ws.Range("L20").Value = 2315
$ws.Range("M20").Value = -763.5833
$ws.Range("N20").Value = -2809
$ws.Range("H64").Value = 644901.06
$ws.Range("I64").Value = 1473224.9
$ws.Range("J64").Value = 649.2222
$ws.Range("K64").Value = 1473224.9
$ws.Range("L64").Value = 649.2222
$ws.Range("M64").Value = -1472999.9
$ws.Range("N64").Value = -1099.2222
$ws.Range("H67").Value = 644901.06
$ws.Range("I67").Value = 1473224.9
$ws.Range("J67").Value = 649.2222
$ws.Range("K67").Value = 1473224.9
$ws.Range("L67").Value = 649.2222
$ws.Range("M67").Value = -1472444.9
$ws.Range("N67").Value = -2209.2222
$ws.Range("H80").Value = 755.2
$ws.Range("I80").Value = 584.6
$ws.Range("K80").Value = 584.6
$ws.Range("M80").Value = 413.4
$ws.Range("H83").Value = 755.2
$ws.Range("I83").Value = 584.6
$ws.Range("K83").Value = 2923
$ws.Range("M83").Value = 2069
$ws.Range("H107").Value = 1268.5555
$ws.Range("I107").Value = 997.1
$ws.Range("J107").Value = 1607.875
$ws.Range("K107").Value = 997.1
$ws.Range("L107").Value = 1607.875
$ws.Range("M107").Value = 922.9
$ws.Range("N107").Value = -5447.875
$ws.Range("H134").Value = 1247.6364
$ws.Range("I134").Value = 1126.7894
$ws.Range("J134").Value = 2013
$ws.Range("K134").Value = 3380.3682
$ws.Range("L134").Value = 6039
$ws.Range("M134").Value = -845.3681999999999
$ws.Range("N134").Value = -11109
$ws.Range("H137").Value = 40856
$ws.Range("J137").Value = 40856
$ws.Range("L137").Value = 40856
$ws.Range("N137").Value = -51056

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43530.785
$ws.Range("I31").Value = 75779.64
$ws.Range("J31").Value = 11281.929
$ws.Range("K31").Value = 75779.64
$ws.Range("L31").Value = 11281.929
$ws.Range("M31").Value = -75484.64
$ws.Range("N31").Value = -11871.929
$ws.Range("H34").Value = 43530.785
$ws.Range("I34").Value = 75779.64
$ws.Range("J34").Value = 11281.929
$ws.Range("K34").Value = 75779.64
$ws.Range("L34").Value = 11281.929
$ws.Range("M34").Value = -75577.64
$ws.Range("N34").Value = -11685.929

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1915.6086
$ws.Range("I102").Value = 1741.4
$ws.Range("K102").Value = 1741.4
$ws.Range("M102").Value = -119.4000000000001
$ws.Range("H113").Value = 1348.9166
$ws.Range("I113").Value = 1036.8334
$ws.Range("J113").Value = 1661
$ws.Range("K113").Value = 1036.8334
$ws.Range("L113").Value = 1661
$ws.Range("M113").Value = 1133.1666
$ws.Range("N113").Value = -6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3569.15
$ws.Range("I122").Value = 3671.2727
$ws.Range("J122").Value = 3444.3333
$ws.Range("K122").Value = 11013.8181
$ws.Range("L122").Value = 10332.9999
$ws.Range("M122").Value = -8563.8181
$ws.Range("N122").Value = -15232.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 640.0769
$ws.Range("I81").Value = 642.1
$ws.Range("J81").Value = 633.3333
$ws.Range("K81").Value = 1284.2
$ws.Range("L81").Value = 1266.6666
$ws.Range("M81").Value = -223.2
$ws.Range("N81").Value = -3388.6666
$ws.Range("H84").Value = 640.0769
$ws.Range("I84").Value = 642.1
$ws.Range("J84").Value = 633.3333
$ws.Range("K84").Value = 6421
$ws.Range("L84").Value = 6333.333000000001
$ws.Range("M84").Value = -1117
$ws.Range("N84").Value = -16941.333
$ws.Range("H132").Value = 1608.8125
$ws.Range("I132").Value = 1386.125
$ws.Range("J132").Value = 1831.5
$ws.Range("K132").Value = 4158.375
$ws.Range("L132").Value = 5494.5
$ws.Range("M132").Value = -1628.375
$ws.Range("N132").Value = -10554.5
$ws.Range("H136").Value = 582.5172
$ws.Range("I136").Value = 271.57144
$ws.Range("J136").Value = 1398.75
$ws.Range("K136").Value = 814.71432
$ws.Range("L136").Value = 4196.25
$ws.Range("M136").Value = 1735.28568
$ws.Range("N136").Value = -9296.25
